# Update the cluster-size data table with refreshed simulation values
# (rows 2-7, columns G/H/N/O/P/Q/R/S/T) and move the sheet selection,
# matching the "added embedded graph to master" re-export of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("N2").Value = 2764.6
$ws.Range("O2").Value = 254.15
$ws.Range("P2").Value = 194.28
$ws.Range("Q2").Value = 192.55
$ws.Range("R2").Value = 193.42
$ws.Range("S2").Value = 2834.57
$ws.Range("T2").Value = 280.27

# --- Row 3 ---
$ws.Range("G3").Value = 2825.59
$ws.Range("H3").Value = 271.02
$ws.Range("N3").Value = 2764.6
$ws.Range("O3").Value = 254.04
$ws.Range("P3").Value = 194.15
$ws.Range("Q3").Value = 192.46
$ws.Range("R3").Value = 193.31
$ws.Range("S3").Value = 2834.57
$ws.Range("T3").Value = 280.273

# --- Row 4 ---
$ws.Range("G4").Value = 2825.86
$ws.Range("H4").Value = 271.29
$ws.Range("N4").Value = 2764.6
$ws.Range("O4").Value = 254.07
$ws.Range("P4").Value = 194.23
$ws.Range("Q4").Value = 192.52
$ws.Range("R4").Value = 193.38
$ws.Range("S4").Value = 2834.6
$ws.Range("T4").Value = 280.3

# --- Row 5 ---
$ws.Range("G5").Value = 2825.9
$ws.Range("H5").Value = 271.3
$ws.Range("N5").Value = 2764.6
$ws.Range("O5").Value = 254.1
$ws.Range("P5").Value = 194.25
$ws.Range("Q5").Value = 192.54
$ws.Range("R5").Value = 193.39
$ws.Range("S5").Value = 2834.6
$ws.Range("T5").Value = 280.3

# --- Row 6 ---
$ws.Range("G6").Value = 2825.86
$ws.Range("H6").Value = 271.29
$ws.Range("N6").Value = 2764.36
$ws.Range("O6").Value = 253.85
$ws.Range("P6").Value = 194.05
$ws.Range("Q6").Value = 192.31
$ws.Range("R6").Value = 193.17
$ws.Range("S6").Value = 2834.56
$ws.Range("T6").Value = 280.25

# --- Row 7 ---
$ws.Range("G7").Value = 2825.9
$ws.Range("H7").Value = 271.3
$ws.Range("N7").Value = 2764.36
$ws.Range("O7").Value = 253.88
$ws.Range("P7").Value = 194.08
$ws.Range("Q7").Value = 192.35
$ws.Range("R7").Value = 193.2
$ws.Range("S7").Value = 2834.57
$ws.Range("T7").Value = 280.27

# Move the visible selection to H7 (matches the saved cursor position in
# the updated workbook).
$ws.Range("H7").Select()
